# "Prefer to pass Evolutions a starter Magic rather than specifying an order value"
#
# Adds a new worksheet named "8" (an 8x8 magic square summing to 260) before
# the existing "9" worksheet, complete with row/column/diagonal sum formulas
# and the same red/green "sum check" conditional formatting used on the
# existing "9" (369) and "30" (13515) magic-square sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: Excel's Font.Color / Interior.Color take a BGR-packed long (same as
# the classic VBA RGB() macro), not a straight 0xRRGGBB value - convert.
# ---------------------------------------------------------------------------
function RgbToBgr($rgbHex) {
    $r = [Convert]::ToInt32($rgbHex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($rgbHex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($rgbHex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$greenFont = RgbToBgr "006100"
$greenFill = RgbToBgr "C6EFCE"
$redFont   = RgbToBgr "9C0006"
$redFill   = RgbToBgr "FFC7CE"

# xlCellValue / xlEqual / xlNotEqual constants
$xlCellValue = 1
$xlEqual = 3
$xlNotEqual = 4

# ---------------------------------------------------------------------------
# Create the new "8" sheet immediately before the "9" sheet.
# ---------------------------------------------------------------------------
$sheet9 = $wb.Worksheets.Item("9")
$ws = $wb.Worksheets.Add($sheet9)
$ws.Name = "8"

# 8x8 magic square (values 1-64), row sums / column sums / both diagonals = 260
$data = @(
  @(12,31,60,33,30,35,39,20),
  @(62,32,16,38,21,7,44,40),
  @(23,28,37,43,11,51,9,58),
  @(25,18,34,48,8,27,54,46),
  @(41,1,45,56,49,6,36,26),
  @(4,53,2,17,63,55,61,5),
  @(64,50,24,15,19,22,14,52),
  @(29,47,42,10,59,57,3,13)
)

for ($r = 0; $r -lt 8; $r++) {
    for ($c = 0; $c -lt 8; $c++) {
        $ws.Cells.Item($r + 1, $c + 2).Value = $data[$r][$c]
    }
}

# Column J: row sums (B:I) for rows 1-8
for ($r = 1; $r -le 8; $r++) {
    $ws.Cells.Item($r, 10).Formula = "=SUM(B${r}:I${r})"
}

# Row 9: column sums (1-8) for columns B-I
for ($c = 2; $c -le 9; $c++) {
    $colLetter = [char](64 + $c)
    $ws.Cells.Item(9, $c).Formula = "=SUM(${colLetter}1:${colLetter}8)"
}

# A9 / J9: the two diagonal sums
$ws.Range("A9").Formula = "=SUM(B8,C7,D6,E5,F4,G3,H2,I1)"
$ws.Range("J9").Formula = "=SUM(B1,C2,D3,E4,F5,G6,H7,I8)"

# ---------------------------------------------------------------------------
# Conditional formatting on the new sheet: green when the sum check equals
# 260, red when it doesn't - same convention as the "9" and "30" sheets.
# ---------------------------------------------------------------------------
$rngRow9 = $ws.Range("A9:J9")
$fcEq = $rngRow9.FormatConditions.Add($xlCellValue, $xlEqual, "=260")
$fcEq.Font.Color = $greenFont
$fcEq.Interior.Color = $greenFill
$fcNe = $rngRow9.FormatConditions.Add($xlCellValue, $xlNotEqual, "=260")
$fcNe.Font.Color = $redFont
$fcNe.Interior.Color = $redFill

$rngColJ = $ws.Range("J1:J8")
$fcEq2 = $rngColJ.FormatConditions.Add($xlCellValue, $xlEqual, "=260")
$fcEq2.Font.Color = $greenFont
$fcEq2.Interior.Color = $greenFill
$fcNe2 = $rngColJ.FormatConditions.Add($xlCellValue, $xlNotEqual, "=260")
$fcNe2.Font.Color = $redFont
$fcNe2.Interior.Color = $redFill

# ---------------------------------------------------------------------------
# Sheet views: the new "8" sheet becomes the active tab (selection D1); the
# "30" sheet keeps its own B31 selection but is no longer the active tab.
# ---------------------------------------------------------------------------
$sheet30 = $wb.Worksheets.Item("30")
$sheet30.Activate()
$sheet30.Range("B31").Select()

$ws.Activate()
$ws.Range("D1").Select()
